$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '34.714.86'
$ws.Range("E2").Value = '  +0.63%  '
$ws.Range("D3").Value = '1.820.97'
$ws.Range("E3").Value = '  +1.49%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '228.85'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.90%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.579'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +3.99%  '
$ws.Range("E7").Value = '  +0.11%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '35.04'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +7.42%  '
$ws.Range("E9").Value = '  +1.67%  '
$ws.Range("E10").Value = '  +0.64%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0952'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.26%  '
$ws.Range("D12").Value = '2.083.79'
$ws.Range("E12").Value = '  +1.48%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '11.41'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +3.32%  '
$ws.Range("D14").Value = '1.836.57'
$ws.Range("E14").Value = '  +2.64%  '
$ws.Range("E15").Value = '  +1.93%  '
$ws.Range("D16").Value = '34.683.85'
$ws.Range("E16").Value = '  +0.64%  '
$ws.Range("E17").Value = '  +1.78%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '69.38'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.71%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '246.93'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.22%  '
$ws.Range("E20").Value = '  +0.09%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '11.61'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +3.43%  '
$ws.Range("E22").Value = '  +0.12%  '
$ws.Range("E23").Value = '  +0.41%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '174.21'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +6.10%  '
$ws.Range("E25").Value = '  +1.65%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '7.56'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +4.11%  '
$ws.Range("E27").Value = '  +2.05%  '
$ws.Range("E28").Value = '  +2.21%  '
$ws.Range("E29").Value = '  -0.08%  '
$ws.Range("E30").Value = '  +2.80%  '
$ws.Range("E31").Value = '  +1.56%  '
$ws.Range("E32").Value = '  +0.81%  '
$ws.Range("E33").Value = '  +0.95%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.85'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +1.39%  '

# Row 35/36: Maker and RenderToken swapped ranking positions
$ws.Range("B35").Value = 'Maker'
$ws.Range("C35").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.402.13'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -2.45%  '
$ws.Range("B36").Value = 'RenderToken'
$ws.Range("C36").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.55'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -1.04%  '

$ws.Range("E37").Value = '  +1.86%  '
$ws.Range("E38").Value = '  -0.81%  '
$ws.Range("E39").Value = '  +0.20%  '
$ws.Range("E40").Value = '  +4.93%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '83.02'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -1.85%  '
$ws.Range("E42").Value = '  +1.82%  '
$ws.Range("E43").Value = '  +0.18%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '13.85'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +2.61%  '
$ws.Range("E45").Value = '  +2.68%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0516'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -2.12%  '
$ws.Range("E47").Value = '  -1.24%  '
$ws.Range("D48").Value = '1.983.53'
$ws.Range("E48").Value = '  +1.71%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '105.22'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.53%  '
$ws.Range("E50").Value = '  -1.01%  '
$ws.Range("E51").Value = '  +0.07%  '
